$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking strings (e.g. "1.00", "7.00")
# are preserved verbatim instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "57.902.58"
$ws.Range("E2").Value = "  +2.02%  "
$ws.Range("D3").Value = "3.048.05"
$ws.Range("E3").Value = "  +2.14%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "525.81"
$ws.Range("E5").Value = "  +5.64%  "
$ws.Range("D6").Value = "142.49"
$ws.Range("E6").Value = "  +5.42%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "0.444"
$ws.Range("E8").Value = "  +4.83%  "
$ws.Range("D9").Value = "7.59"
$ws.Range("E9").Value = "  +3.81%  "
$ws.Range("D10").Value = "0.113"
$ws.Range("E10").Value = "  +8.76%  "
$ws.Range("D11").Value = "0.369"
$ws.Range("E11").Value = "  +5.27%  "
$ws.Range("E12").Value = "  +2.68%  "
$ws.Range("D13").Value = "3.571.85"
$ws.Range("E13").Value = "  +2.39%  "
$ws.Range("D14").Value = "26.75"
$ws.Range("E14").Value = "  +7.67%  "
$ws.Range("E15").Value = "  +17.71%  "
$ws.Range("D16").Value = "57.870.06"
$ws.Range("E16").Value = "  +2.15%  "
$ws.Range("D17").Value = "6.20"
$ws.Range("E17").Value = "  +6.90%  "
$ws.Range("D18").Value = "3.046.84"
$ws.Range("E18").Value = "  +2.32%  "
$ws.Range("D19").Value = "12.94"
$ws.Range("E19").Value = "  +5.16%  "
$ws.Range("D20").Value = "8.20"
$ws.Range("E20").Value = "  +6.11%  "
$ws.Range("D21").Value = "340.91"
$ws.Range("E21").Value = "  +5.28%  "
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Value = "0.499"
$ws.Range("E23").Value = "  +7.99%  "
$ws.Range("D24").Value = "65.37"
$ws.Range("E24").Value = "  +6.78%  "
$ws.Range("E25").Value = "  +7.07%  "
$ws.Range("D26").Value = "0.0₃0968"
$ws.Range("E26").Value = "  +7.42%  "
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").Value = "7.00"
$ws.Range("E28").Value = "  +7.55%  "
$ws.Range("D29").Value = "7.33"
$ws.Range("E29").Value = "  +8.02%  "
$ws.Range("E30").Value = "  +7.50%  "
$ws.Range("E31").Value = "  +5.69%  "
$ws.Range("D32").Value = "21.07"
$ws.Range("E32").Value = "  +6.00%  "
$ws.Range("D33").Value = "156.74"
$ws.Range("E33").Value = "  +1.09%  "
$ws.Range("D34").Value = "4.73"
$ws.Range("E34").Value = "  +5.96%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").Value = "5.91"
$ws.Range("E35").Value = "  +5.25%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "1.34"
$ws.Range("E36").Value = "  +4.39%  "
$ws.Range("D37").Value = "25.92"
$ws.Range("E37").Value = "  +10.44%  "
$ws.Range("D38").Value = "0.0690"
$ws.Range("E38").Value = "  +2.68%  "
$ws.Range("D39").Value = "3.082.23"
$ws.Range("E39").Value = "  +2.22%  "
$ws.Range("D40").Value = "37.72"
$ws.Range("E40").Value = "  +1.72%  "
$ws.Range("D41").Value = "3.86"
$ws.Range("E41").Value = "  +8.50%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("E43").Value = "  +4.60%  "
$ws.Range("D44").Value = "0.663"
$ws.Range("E44").Value = "  +4.32%  "
$ws.Range("D45").Value = "2.325.71"
$ws.Range("E45").Value = "  +5.21%  "
$ws.Range("E46").Value = "  +1.66%  "
$ws.Range("D47").Value = "2.04"
$ws.Range("E47").Value = "  +3.47%  "
$ws.Range("E48").Value = "  +4.97%  "
$ws.Range("D49").Value = "6.05"
$ws.Range("E49").Value = "  +6.06%  "
$ws.Range("D50").Value = "20.03"
$ws.Range("E50").Value = "  +4.29%  "
$ws.Range("D51").Value = "0.0894"
$ws.Range("E51").Value = "  +5.62%  "

# Restore column D to the default (unstyled) cell style now that the text values are set,
# so the cells do not retain an explicit style reference.
$ws.Range("D2:D51").Style = "Normal"

